$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $c = $ws.Range($Range)
    $c.NumberFormat = "@"
    $c.Value = $Value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '60.119.35'
Set-TextValue 'E2' '  -0.49%  '
Set-TextValue 'D3' '2.607.05'
Set-TextValue 'E3' '  +0.71%  '
Set-TextValue 'D5' '519.47'
Set-TextValue 'E5' '  +0.95%  '
Set-TextValue 'D6' '148.62'
Set-TextValue 'E6' '  -3.09%  '
Set-TextValue 'E7' '  +0.07%  '
Set-TextValue 'E8' '  -4.78%  '
Set-TextValue 'D9' '2.614.20'
Set-TextValue 'E9' '  +0.82%  '
Set-TextValue 'D10' '6.28'
Set-TextValue 'E10' '  -4.95%  '
Set-TextValue 'E11' '  +0.24%  '
Set-TextValue 'E12' '  -1.54%  '
Set-TextValue 'D13' '0.129'
Set-TextValue 'E13' '  -0.68%  '
Set-TextValue 'D14' '3.064.98'
Set-TextValue 'E14' '  +0.69%  '
Set-TextValue 'D15' '60.182.48'
Set-TextValue 'E15' '  -0.37%  '
Set-TextValue 'D16' '21.16'
Set-TextValue 'E16' '  -2.51%  '
Set-TextValue 'D17' '0.0000137'
Set-TextValue 'E17' '  -1.71%  '
Set-TextValue 'D18' '2.600.63'
Set-TextValue 'E18' '  +0.30%  '
Set-TextValue 'D19' '4.62'
Set-TextValue 'E19' '  -2.76%  '
Set-TextValue 'D20' '341.92'
Set-TextValue 'E20' '  -2.75%  '
Set-TextValue 'D21' '10.38'
Set-TextValue 'E21' '  -1.72%  '
Set-TextValue 'E22' '  -1.93%  '
Set-TextValue 'D23' '0.997'
Set-TextValue 'E23' '  -0.15%  '
Set-TextValue 'D24' '60.25'
Set-TextValue 'E24' '  -1.05%  '
Set-TextValue 'E25' '  -2.27%  '
Set-TextValue 'D26' '0.998'
Set-TextValue 'E26' '  +0.54%  '
Set-TextValue 'D27' '0.161'
Set-TextValue 'E27' '  -2.74%  '
Set-TextValue 'E28' '  -3.84%  '
Set-TextValue 'D29' '7.05'
Set-TextValue 'E29' '  -3.14%  '
Set-TextValue 'E30' '  -0.02%  '
Set-TextValue 'D31' '6.01'
Set-TextValue 'E31' '  -3.75%  '
Set-TextValue 'D32' '1.58'
Set-TextValue 'E32' '  -0.20%  '
Set-TextValue 'E33' '  -2.59%  '
Set-TextValue 'D34' '149.70'
Set-TextValue 'E34' '  -0.44%  '
Set-TextValue 'D35' '3.96'
Set-TextValue 'E35' '  -2.60%  '
Set-TextValue 'D36' '0.905'
Set-TextValue 'E36' '  -3.13%  '
Set-TextValue 'D37' '1.13'
Set-TextValue 'E37' '  -4.75%  '
Set-TextValue 'D38' '0.860'
Set-TextValue 'E38' '  +3.01%  '
Set-TextValue 'D39' '36.40'
Set-TextValue 'E39' '  +0.34%  '
Set-TextValue 'D40' '1.43'
Set-TextValue 'E40' '  -3.54%  '
Set-TextValue 'E41' '  -3.97%  '
Set-TextValue 'D42' '287.07'
Set-TextValue 'E42' '  +1.07%  '
Set-TextValue 'D43' '0.624'
Set-TextValue 'E43' '  +0.29%  '
Set-TextValue 'E44' '  -1.52%  '
Set-TextValue 'E45' '  +0.21%  '
Set-TextValue 'D46' '0.0546'
Set-TextValue 'E46' '  -2.19%  '
Set-TextValue 'D47' '19.44'
Set-TextValue 'E47' '  +0.20%  '
Set-TextValue 'B48' 'WhiteBITCoin'
Set-TextValue 'C48' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D48' '10.39'
Set-TextValue 'E48' '  +0.86%  '
Set-TextValue 'B49' 'VeChain'
Set-TextValue 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D49' '0.0230'
Set-TextValue 'E49' '  -2.42%  '
Set-TextValue 'D50' '4.64'
Set-TextValue 'E50' '  -4.10%  '
Set-TextValue 'D51' '1.952.07'
Set-TextValue 'E51' '  -0.84%  '
